$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.565.44'
$ws.Range("E2").Value = '  +1.76%  '
$ws.Range("D3").Value = '2.215.33'
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.619'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.88'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.06%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.594'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.27'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0927'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.88'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E14").Value = '  -2.19%  '
$ws.Range("D15").Value = '2.548.52'
$ws.Range("E15").Value = '  -0.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.62'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.06%  '
$ws.Range("D17").Value = '2.218.12'
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.801'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.40%  '
$ws.Range("D19").Value = '42.415.51'
$ws.Range("E19").Value = '  +1.66%  '
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.70%  '
$ws.Range("E22").Value = '  -3.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -8.21%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '229.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.73%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.07%  '
$ws.Range("E29").Value = '  -1.96%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '172.53'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +19.94%  '
$ws.Range("E32").Value = '  -4.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0794'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.67%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.108'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.55%  '
$ws.Range("E38").Value = '  +3.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0321'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.49'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '60.29'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.63%  '
$ws.Range("E44").Value = '  +0.91%  '
$ws.Range("E45").Value = '  -0.86%  '
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '99.29'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.11'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.427'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +16.67%  '
